# "abis pre processing, progress cek catetan"
#
# Changes applied to the "jumlah missing values" sheet:
#  1. Rows 2-26 (the data rows whose MissingValues count the pre-processing
#     pass already cleaned up) get a slightly taller row height, from 18.75
#     to 19.5 points - matching what Excel writes after re-touching those
#     rows during the cleanup pass.
#  2. The MissingValues numbers in column B (B2:B75) get their font color
#     pinned to explicit black (RGB 000000) instead of the theme-1 color,
#     which is what Excel does when it re-saves and collapses the
#     (accidental) duplicate font definition that used to back those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the row height for rows 2 through 26.
$ws.Range("A2:A26").EntireRow.RowHeight = 19.5

# 2. Make the numeric MissingValues column use an explicit black font color.
$ws.Range("B2:B75").Font.Color = 0
